$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Remove "Juan Bernal Jimenez" from the FSR list
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $v = $ws.Cells.Item($r, 1).Value()
    if ($v -eq "Juan Bernal Jimenez") {
        $ws.Rows($r).Delete()
        break
    }
}

# Add the new FSR names to the table
$newNames = @(
    "Ana Cristina Soto Ruvalcaba",
    "Carlos Vargas",
    "Eduardo  Artero",
    "Gema Llorente",
    "Luis Moran",
    "Manuela Martin",
    "Maria Dolores Redondo",
    "Marta Sanchez",
    "Nicolás Capelastegui Rojo",
    "Susana Gonzalez",
    "Asier Garcia Doncel",
    "Lara Morote Moreno",
    "Jose Antonio Gragera Cano",
    "Julio Jesus Martinez Romero"
)
foreach ($name in $newNames) {
    $newRow = $lo.ListRows.Add()
    $newRow.Range.Value = $name
}

# Keep the table sorted alphabetically (A-Z) on the FSR column, as before
$lo.DataBodyRange.Sort($lo.DataBodyRange.Columns.Item(1))

# Highlight the newly added / noteworthy entries in yellow
$highlightNames = @("Ana Cristina Soto Ruvalcaba", "Julio Jesus Martinez Romero", "Monica Gouveia")
for ($i = 1; $i -le $lo.ListRows.Count; $i++) {
    $rowRange = $lo.ListRows.Item($i).Range
    if ($highlightNames -contains $rowRange.Value()) {
        $rowRange.Interior.Color = 65535
    }
}

# Widen column A to fit the longer names
$ws.Columns.Item(1).ColumnWidth = 26.75
